# Auto-generated market-data refresh script
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N) for the
# affected leve rows on each class sheet, as produced by the scheduled price-scrape runner.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1726.2222
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 1567
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 1567
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -2219
$ws.Range("H64").Value = 1976181.9
$ws.Range("I64").Value = 2567820.5
$ws.Range("J64").Value = 4053.3333
$ws.Range("K64").Value = 2567820.5
$ws.Range("L64").Value = 4053.3333
$ws.Range("M64").Value = -2567572.5
$ws.Range("N64").Value = -4549.3333
$ws.Range("H67").Value = 1976181.9
$ws.Range("I67").Value = 2567820.5
$ws.Range("J67").Value = 4053.3333
$ws.Range("K67").Value = 2567820.5
$ws.Range("L67").Value = 4053.3333
$ws.Range("M67").Value = -2566962.5
$ws.Range("N67").Value = -5769.3333
$ws.Range("H76").Value = 3979.4443
$ws.Range("I76").Value = 3830.7144
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 3830.7144
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -3515.7144
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 3979.4443
$ws.Range("I79").Value = 3830.7144
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 3830.7144
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -2738.7144
$ws.Range("N79").Value = -6684
$ws.Range("H141").Value = 1216.0714
$ws.Range("I141").Value = 1078.8462
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 3236.5386
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 1943.4614
$ws.Range("N141").Value = -19360

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 50
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 50
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = None
$ws.Range("H32").Value = 3032.22
$ws.Range("I32").Value = 3032.22
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3032.22
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = None
$ws.Range("H45").Value = 2274239.5
$ws.Range("I45").Value = 2526710.8
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2526710.8
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -2526333.8
$ws.Range("N45").Value = -2754
$ws.Range("H110").Value = 912.2174
$ws.Range("I110").Value = 634.05
$ws.Range("J110").Value = 2766.6667
$ws.Range("K110").Value = 634.05
$ws.Range("L110").Value = 2766.6667
$ws.Range("M110").Value = 1410.95
$ws.Range("N110").Value = -6856.6667
$ws.Range("N25").ClearContents()
$ws.Range("N32").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1880.6666
$ws.Range("I3").Value = 1800
$ws.Range("J3").Value = 1945.2
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 1945.2
$ws.Range("M3").Value = -1686
$ws.Range("N3").Value = -2173.2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5228
$ws.Range("H54").Value = 4195.385
$ws.Range("I54").Value = 2854
$ws.Range("J54").Value = 8666.666999999999
$ws.Range("K54").Value = 2854
$ws.Range("L54").Value = 8666.666999999999
$ws.Range("M54").Value = -2370
$ws.Range("N54").Value = -9634.666999999999
$ws.Range("H134").Value = 4864.472
$ws.Range("I134").Value = 3970.2856
$ws.Range("J134").Value = 5433.5
$ws.Range("K134").Value = 11910.8568
$ws.Range("L134").Value = 16300.5
$ws.Range("M134").Value = -9375.856800000001
$ws.Range("N134").Value = -21370.5
$ws.Range("M29").ClearContents()

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 5500
$ws.Range("J29").Value = 5500
$ws.Range("L29").Value = 5500
$ws.Range("N29").Value = -6086
$ws.Range("H54").Value = 29962.666
$ws.Range("J54").Value = 29962.666
$ws.Range("L54").Value = 29962.666
$ws.Range("N54").Value = -31278.666
$ws.Range("H62").Value = 2533.6365
$ws.Range("I62").Value = 2195
$ws.Range("J62").Value = 2727.1428
$ws.Range("K62").Value = 2195
$ws.Range("L62").Value = 2727.1428
$ws.Range("M62").Value = -1571
$ws.Range("N62").Value = -3975.1428
$ws.Range("H65").Value = 2533.6365
$ws.Range("I65").Value = 2195
$ws.Range("J65").Value = 2727.1428
$ws.Range("K65").Value = 10975
$ws.Range("L65").Value = 13635.714
$ws.Range("M65").Value = -7855
$ws.Range("N65").Value = -19875.714
$ws.Range("H70").Value = 30025.715
$ws.Range("J70").Value = 30025.715
$ws.Range("L70").Value = 30025.715
$ws.Range("N70").Value = -30655.715
$ws.Range("H73").Value = 30025.715
$ws.Range("J73").Value = 30025.715
$ws.Range("L73").Value = 30025.715
$ws.Range("N73").Value = -32209.715
$ws.Range("H76").Value = 13276
$ws.Range("I76").Value = 13276
$ws.Range("K76").Value = 13276
$ws.Range("M76").Value = -12961
$ws.Range("H79").Value = 13276
$ws.Range("I79").Value = 13276
$ws.Range("K79").Value = 13276
$ws.Range("M79").Value = -12184
$ws.Range("H99").Value = 1994.7391
$ws.Range("I99").Value = 1761.091
$ws.Range("J99").Value = 2208.9167
$ws.Range("K99").Value = 1761.091
$ws.Range("L99").Value = 2208.9167
$ws.Range("M99").Value = -263.0909999999999
$ws.Range("N99").Value = -5204.9167
$ws.Range("H126").Value = 1994.7391
$ws.Range("I126").Value = 1761.091
$ws.Range("J126").Value = 2208.9167
$ws.Range("K126").Value = 5283.272999999999
$ws.Range("L126").Value = 6626.750100000001
$ws.Range("M126").Value = -2813.272999999999
$ws.Range("N126").Value = -11566.7501
$ws.Range("H134").Value = 662795.9
$ws.Range("I134").Value = 1543.25
$ws.Range("J134").Value = 5952817
$ws.Range("K134").Value = 4629.75
$ws.Range("L134").Value = 17858451
$ws.Range("M134").Value = -2094.75
$ws.Range("N134").Value = -17863521

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 37037810
$ws.Range("I113").Value = 125000540
$ws.Range("J113").Value = 865.7895
$ws.Range("K113").Value = 375001620
$ws.Range("L113").Value = 2597.3685
$ws.Range("M113").Value = -374999450
$ws.Range("N113").Value = -6937.3685

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 6710
$ws.Range("J24").Value = 6710
$ws.Range("L24").Value = 6710
$ws.Range("N24").Value = -7056
$ws.Range("H70").Value = 14208.421
$ws.Range("I70").Value = 42800
$ws.Range("K70").Value = 42800
$ws.Range("M70").Value = -42530
$ws.Range("H73").Value = 14208.421
$ws.Range("I73").Value = 42800
$ws.Range("K73").Value = 42800
$ws.Range("M73").Value = -41864
$ws.Range("H102").Value = 2586.3635
$ws.Range("I102").Value = 3163.238
$ws.Range("J102").Value = 1576.8334
$ws.Range("K102").Value = 3163.238
$ws.Range("L102").Value = 1576.8334
$ws.Range("M102").Value = -1541.238
$ws.Range("N102").Value = -4820.8334

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1267.8108
$ws.Range("I16").Value = 1222.2333
$ws.Range("J16").Value = 1463.1428
$ws.Range("K16").Value = 1222.2333
$ws.Range("L16").Value = 1463.1428
$ws.Range("M16").Value = -1052.2333
$ws.Range("N16").Value = -1803.1428
$ws.Range("H22").Value = 1153.5294
$ws.Range("I22").Value = 574
$ws.Range("K22").Value = 574
$ws.Range("M22").Value = -279
$ws.Range("H27").Value = 1153.5294
$ws.Range("I27").Value = 574
$ws.Range("K27").Value = 574
$ws.Range("M27").Value = -467
$ws.Range("H54").Value = 23472
$ws.Range("J54").Value = 23472
$ws.Range("L54").Value = 23472
$ws.Range("N54").Value = -24760
$ws.Range("H68").Value = 1575.6666
$ws.Range("I68").Value = 1701.1538
$ws.Range("J68").Value = 1427.3636
$ws.Range("K68").Value = 1701.1538
$ws.Range("L68").Value = 1427.3636
$ws.Range("M68").Value = -952.1538
$ws.Range("N68").Value = -2925.3636
$ws.Range("H71").Value = 1575.6666
$ws.Range("I71").Value = 1701.1538
$ws.Range("J71").Value = 1427.3636
$ws.Range("K71").Value = 8505.769
$ws.Range("L71").Value = 7136.817999999999
$ws.Range("M71").Value = -4761.769
$ws.Range("N71").Value = -14624.818
$ws.Range("H132").Value = 12511632
$ws.Range("I132").Value = 9731.925999999999
$ws.Range("J132").Value = 38477116
$ws.Range("K132").Value = 29195.778
$ws.Range("L132").Value = 115431348
$ws.Range("M132").Value = -26665.778
$ws.Range("N132").Value = -115436408

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5538.1177
$ws.Range("I62").Value = 5475
$ws.Range("J62").Value = 5594.222
$ws.Range("K62").Value = 5475
$ws.Range("L62").Value = 5594.222
$ws.Range("M62").Value = -4851
$ws.Range("N62").Value = -6842.222
$ws.Range("H65").Value = 5538.1177
$ws.Range("I65").Value = 5475
$ws.Range("J65").Value = 5594.222
$ws.Range("K65").Value = 27375
$ws.Range("L65").Value = 27971.11
$ws.Range("M65").Value = -24255
$ws.Range("N65").Value = -34211.11
$ws.Range("H70").Value = 22000
$ws.Range("J70").Value = 22000
$ws.Range("L70").Value = 22000
$ws.Range("N70").Value = -22630
$ws.Range("H73").Value = 22000
$ws.Range("J73").Value = 22000
$ws.Range("L73").Value = 22000
$ws.Range("N73").Value = -24184
$ws.Range("H132").Value = 1351.7537
$ws.Range("I132").Value = 1166.2909
$ws.Range("J132").Value = 2080.3572
$ws.Range("K132").Value = 3498.8727
$ws.Range("L132").Value = 6241.071599999999
$ws.Range("M132").Value = -968.8726999999999
$ws.Range("N132").Value = -11301.0716
